# Update 227C class list
# - Insert a new roster row (Peijin Gao) before the "Andres" row, shifting the
#   existing rows down by one.
# - Paste the raw registrar export (ID / Name / Email / Dept / Level / Grade)
#   into column E for every student row.
# - Apply the small "Arial Unicode MS" 10pt font used for that pasted column.
# - Bump the row height slightly (17 -> 18) to fit the new column's text.
# - Restore the print orientation and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Shift the existing "Andres" .. "Anthony" rows (5-18) down to (6-19) so a
#    new row opens up at row 5 for the not-yet-rostered student (Peijin).
#    Walk bottom-up so we never clobber a row before it has been read.
# ---------------------------------------------------------------------------
for ($r = 18; $r -ge 5; $r--) {
    $dest = $r + 1
    $ws.Cells.Item($dest, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Copy()
    $ws.Cells.Item($dest, 2).PasteSpecial(-4122)
    $ws.Cells.Item($dest, 2).Value = $ws.Cells.Item($r, 2).Value2
}

# Row 5 becomes the new row: only column A (name) is populated; there is no
# roster (column B) entry yet for this student. The name itself is filled in
# after the column-E pastes below so the shared-string table ends up in the
# same append order as the authored workbook.
$ws.Cells.Item(5, 2).ClearContents()

# ---------------------------------------------------------------------------
# 2) Paste the raw export string into column E, rows 2-19, and give that
#    column its own small font (first cell sets it up, the rest copy the
#    resulting format so only one new font/style entry is created).
# ---------------------------------------------------------------------------
$rawRows = @(
    "85121744        CHEN, CHAORONG  CHAORONC@UCI.EDU        MCSB    G6      GR",
    "15841424        DE ROBLES, GABRIELA     GDEROBLE@UCI.EDU        MCSB    G5      GR",
    "21244597        ELDEEN, SARAH   SELDEEN@UCI.EDU MCSB    G6      GR",
    "84549672        GAO, PEIJUN     PEIJUNG@UCI.EDU MATH    JR      GR",
    "25167480        GUERRERO RAMIREZ, ANDRES FELIPE AFGUERRE@UCI.EDU        MCSB    G6      GR",
    "18249792        IYER, VIGNESH HARIHARAN VHIYER@UCI.EDU  MATH    G5      GR",
    "49527062        LAI, LULU       LLAI7@UCI.EDU   MCSB    G6      GR",
    "84409046        NGUYEN, THI THU THAO    THAOTN18@UCI.EDU        MCSB    G6      GR",
    "61007128        OLARANONT, NONTHAKORN   NOLARANO@UCI.EDU        MATH    G6      GR",
    "41761116        RAHIMZADEH, NEGIN       NRAHIMZA@UCI.EDU        MCSB    G6      GR",
    "13522870        RUAN, YIBIAO    YIBIAOR@UCI.EDU MATH    G5      GR",
    "90727596        SHAO, WEI       SHAOW6@UCI.EDU  MCSB    G6      GR",
    "19919996        SILKWOOD, KAI   KSILKWOO@UCI.EDU        MCSB    G6      GR",
    "12442257        TAN, PEI        PEIT3@UCI.EDU   MCSB    G5      GR",
    "63083899        TRAN, NHAT THANH VAN    NHATTT@UCI.EDU  MATH    G6      GR",
    "92807178        WANG, TIANHONG  TIANHOW1@UCI.EDU        MATH+BIO SCI    SR      GR",
    "22924311        WEI, SIYANG     SIYANGW2@UCI.EDU        MATH    G6      GR",
    "39367708        ZAMORA, ANTHONY ZAMORAA3@UCI.EDU        MATH    G6      GR"
)

$firstCell = $ws.Cells.Item(2, 5)
$firstCell.Value = $rawRows[0]
$firstCell.Font.Name = "Arial Unicode MS"
$firstCell.Font.Size = 10
$firstCell.Font.Color = 0
$firstCell.Copy()

for ($i = 1; $i -lt $rawRows.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 5)
    $cell.Value = $rawRows[$i]
    $cell.PasteSpecial(-4122)
}

# Now that every other new string has been appended, record the new
# student's short name in column A.
$ws.Cells.Item(5, 1).Value = "Peijin"

# ---------------------------------------------------------------------------
# 3) Bump the row height for the student rows to fit the new column.
# ---------------------------------------------------------------------------
$ws.Range("A2:I19").RowHeight = 18

# ---------------------------------------------------------------------------
# 4) Misc sheet-level touch-ups matching the authored change.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("I26").Select()
